# Week 15 (logged) + Week 16 (simulated) update for Jets Players Data
$wb = $excel.ActiveWorkbook

$rushing = $wb.Worksheets.Item("Rushing")
$receiving = $wb.Worksheets.Item("Receiving")

# Fix player name typo: Mi.Carter -> M.Carter
$rushing.Cells.Item(8, 2).Value = "M.Carter"
$receiving.Cells.Item(4, 2).Value = "M.Carter"

# Rushing sheet data updates
$rushing.Range("D2").Value = 4
$rushing.Range("E2").Value = 3
$rushing.Range("F2").Value = 6

$rushing.Range("C6").Value = 35
$rushing.Range("D6").Value = 23
$rushing.Range("F6").Value = 8

$rushing.Range("C8").Value = 70
$rushing.Range("D8").Value = 50
$rushing.Range("E8").Value = 20
$rushing.Range("F8").Value = 25

$rushing.Range("C10").Value = 8
$rushing.Range("E10").Value = 2

$rushing.Range("C12").Value = 3
$rushing.Range("E12").Value = 1
$rushing.Range("F12").Value = 2

# Receiving sheet data updates
$receiving.Range("C4").Value = 50
$receiving.Range("D4").Value = 45
$receiving.Range("E4").Value = 3
$receiving.Range("F4").Value = 2
$receiving.Range("G4").Value = 1
$receiving.Range("H4").Value = 1

$receiving.Range("C6").Value = 87
$receiving.Range("D6").Value = 65

$receiving.Range("C7").Value = 30
$receiving.Range("D7").Value = 16
$receiving.Range("G7").Value = 7

$receiving.Range("C9").Value = 40
$receiving.Range("D9").Value = 28
$receiving.Range("G9").Value = 5
$receiving.Range("H9").Value = 3

$receiving.Range("C10").Value = 12
$receiving.Range("E10").Value = 7

$receiving.Range("C14").Value = 15
$receiving.Range("D14").Value = 10
$receiving.Range("E14").Value = 3
$receiving.Range("F14").Value = 2

$receiving.Range("C15").Value = 36
$receiving.Range("D15").Value = 23
$receiving.Range("E15").Value = 6
$receiving.Range("F15").Value = 4

$receiving.Range("C16").Value = 5
$receiving.Range("D16").Value = 3

# Active sheet is now Receiving, with selection at C19
$rushing.Range("C26").Select()
$receiving.Select()
$receiving.Range("C19").Select()
